$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-Mfd3")

# Remove the extra, unused sheets ("Лист2" / "Лист3") - delete from the
# end so indices stay valid while deleting.
$excel.DisplayAlerts = $false
for ($i = $wb.Worksheets.Count; $i -ge 1; $i--) {
    $sheet = $wb.Worksheets.Item($i)
    if ($sheet.Name -ne "Bill of Materials-Mfd3") {
        $sheet.Delete()
    }
}
$excel.DisplayAlerts = $true

# The BOM now reports JLCPCB part numbers instead of LCSC ones.
$ws.Range("F1").Value = "JLCPCB Part #"

# TS-H003 tact switch row now has its JLCPCB part number filled in. Set the
# value first, then copy the sibling (already-filled) cell's formatting onto
# it so it matches the "has a part number" styling used by the rest of
# column F instead of the "blank" styling it started with.
$ws.Range("F16").Value = "C2884874"
$ws.Range("F2").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the last interactive selection.
$ws.Range("G9").Select()
